$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new boolean value in A7 (FALSE)
$ws.Range("A7").Value = $false

# Update the active selection to A8 to match the new state
$ws.Range("A8").Select()
